# Apply text edits described by the commit:
# fix: remove fabricated data, clear language, add migration section

$p = $ppt.ActivePresentation

# Slide 10: "10,000+ D&A Professionals" -> "Global Scale, Local Expertise"
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(17).TextFrame.TextRange.Text = "Global Scale, Local Expertise"

# Slide 10: "35+ countries. Deep IBIOL sectoral knowledge. 75% of Fortune Global 100 served."
#   -> "[30B]+ organization. ~200,000 employees across 70+ countries. 75% of Fortune Global 100 served."
#   (dollar sign escaped below so PowerShell treats it literally, not as a variable)
$s10.Shapes.Item(18).TextFrame.TextRange.Text = "`$30B+ organization. ~200,000 employees across 70+ countries. 75% of Fortune Global 100 served."

# Slide 9: "PHASE 1  ·  H1 2026" -> "PHASE 1  ·  1º SEM 2026"
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(6).TextFrame.TextRange.Text = "PHASE 1  ·  1º SEM 2026"

# Slide 9: "PHASE 2  ·  H2 2026" -> "PHASE 2  ·  2º SEM 2026"
$s9.Shapes.Item(20).TextFrame.TextRange.Text = "PHASE 2  ·  2º SEM 2026"
